# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# 股票 (stock) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): copy the existing header style onto the new cells ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-4) ---
# Format the date column as text first so the "yyyy-mm-dd" string isn't
# silently re-interpreted as a date serial number.
$ws.Range("H2:H4").NumberFormat = "@"

$ws.Range("H2").Value = "2012-04-26"
$ws.Range("I2").Value = "黃昭順"
$ws.Range("J2").Value = 665

$ws.Range("H3").Value = "2012-04-26"
$ws.Range("I3").Value = "黃昭順"
$ws.Range("J3").Value = 665

$ws.Range("H4").Value = "2012-04-26"
$ws.Range("I4").Value = "黃昭順"
$ws.Range("J4").Value = 665
